$wb = $excel.ActiveWorkbook

# Hunk 0: sheet ALC, row 86
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 4987.0586
$ws.Range("I86").Value = 3640.889
$ws.Range("J86").Value = 6501.5
$ws.Range("K86").Value = 3640.889
$ws.Range("L86").Value = 6501.5
$ws.Range("M86").Value = -2517.889
$ws.Range("N86").Value = -8747.5

# Hunk 1: sheet ALC, row 89
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 4987.0586
$ws.Range("I89").Value = 3640.889
$ws.Range("J89").Value = 6501.5
$ws.Range("K89").Value = 18204.445
$ws.Range("L89").Value = 32507.5
$ws.Range("M89").Value = -12588.445
$ws.Range("N89").Value = -43739.5

# Hunk 2: sheet ARM, row 6
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 16125.5
$ws.Range("I6").Value = 16429.143
$ws.Range("J6").Value = 14000
$ws.Range("K6").Value = 16429.143
$ws.Range("L6").Value = 14000
$ws.Range("M6").Value = -16256.143
$ws.Range("N6").Value = -14346

# Hunk 3: sheet ARM, row 35
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H35").Value = 5779
$ws.Range("I35").Value = 5779
$ws.Range("K35").Value = 5779
$ws.Range("M35").Value = -5373

# Hunk 4: sheet BSM, row 5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 2000
$ws.Range("I5").Value = 2000
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 2000
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -1887
$ws.Range("N5").ClearContents()

# Hunk 5: sheet BSM, row 37
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H37").Value = 8131.154
$ws.Range("I37").Value = 2468.5715
$ws.Range("J37").Value = 14737.5
$ws.Range("K37").Value = 2468.5715
$ws.Range("L37").Value = 14737.5
$ws.Range("M37").Value = -2331.5715
$ws.Range("N37").Value = -15011.5

# Hunk 6: sheet BSM, row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5260.8945
$ws.Range("I134").Value = 5464.923
$ws.Range("J134").Value = 4818.8335
$ws.Range("K134").Value = 16394.769
$ws.Range("L134").Value = 14456.5005
$ws.Range("M134").Value = -13859.769
$ws.Range("N134").Value = -19526.5005

# Hunk 7: sheet BSM, row 140
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# Hunk 8: sheet BSM, row 141
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

# Hunk 9: sheet CRP, row 12
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 202.5
$ws.Range("I12").Value = 202.5
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 202.5
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -32.5
$ws.Range("N12").ClearContents()

# Hunk 10: sheet CRP, row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1364
$ws.Range("I31").Value = 943.5
$ws.Range("K31").Value = 943.5
$ws.Range("M31").Value = -648.5

# Hunk 11: sheet CRP, row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1364
$ws.Range("I34").Value = 943.5
$ws.Range("K34").Value = 943.5
$ws.Range("M34").Value = -741.5

# Hunk 12: sheet CRP, row 39
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H39").Value = 8825
$ws.Range("I39").Value = 8825
$ws.Range("K39").Value = 8825
$ws.Range("M39").Value = -8434

# Hunk 13: sheet CRP, row 49
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H49").Value = 8825
$ws.Range("I49").Value = 8825
$ws.Range("K49").Value = 8825
$ws.Range("M49").Value = -8643

# Hunk 14: sheet CRP, row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3666.348
$ws.Range("I132").Value = 2860.7058
$ws.Range("J132").Value = 5949
$ws.Range("K132").Value = 8582.117400000001
$ws.Range("L132").Value = 17847
$ws.Range("M132").Value = -6052.117400000001
$ws.Range("N132").Value = -22907

# Hunk 15: sheet CUL, row 21
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H21").Value = 1534.4615
$ws.Range("I21").Value = 2078.8
$ws.Range("J21").Value = 1454.4117
$ws.Range("K21").Value = 6236.400000000001
$ws.Range("L21").Value = 4363.2351
$ws.Range("M21").Value = -6063.400000000001
$ws.Range("N21").Value = -4709.2351

# Hunk 16: sheet CUL, row 22
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 1216.2162
$ws.Range("I22").Value = 780
$ws.Range("J22").Value = 1284.375
$ws.Range("K22").Value = 2340
$ws.Range("L22").Value = 3853.125
$ws.Range("M22").Value = -2171
$ws.Range("N22").Value = -4191.125

# Hunk 17: sheet CUL, row 27
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H27").Value = 1216.2162
$ws.Range("I27").Value = 780
$ws.Range("J27").Value = 1284.375
$ws.Range("K27").Value = 2340
$ws.Range("L27").Value = 3853.125
$ws.Range("M27").Value = -2238
$ws.Range("N27").Value = -4057.125

# Hunk 18: sheet CUL, row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 971.8108
$ws.Range("I131").Value = 439.66666
$ws.Range("J131").Value = 1074.8064
$ws.Range("K131").Value = 1318.99998
$ws.Range("L131").Value = 3224.4192
$ws.Range("M131").Value = 3721.00002
$ws.Range("N131").Value = -13304.4192

# Hunk 19: sheet GSM, row 5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 24869.4
$ws.Range("I5").Value = 90000
$ws.Range("J5").Value = 8586.75
$ws.Range("K5").Value = 90000
$ws.Range("L5").Value = 8586.75
$ws.Range("M5").Value = -89888
$ws.Range("N5").Value = -8810.75

# Hunk 20: sheet GSM, row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5374.6875
$ws.Range("I70").Value = 5227.273
$ws.Range("J70").Value = 5699
$ws.Range("K70").Value = 5227.273
$ws.Range("L70").Value = 5699
$ws.Range("M70").Value = -4957.273
$ws.Range("N70").Value = -6239

# Hunk 21: sheet GSM, row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 5374.6875
$ws.Range("I73").Value = 5227.273
$ws.Range("J73").Value = 5699
$ws.Range("K73").Value = 5227.273
$ws.Range("L73").Value = 5699
$ws.Range("M73").Value = -4291.273
$ws.Range("N73").Value = -7571

# Hunk 22: sheet GSM, row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2511.3225
$ws.Range("I80").Value = 2540.682
$ws.Range("J80").Value = 2439.5557
$ws.Range("K80").Value = 2540.682
$ws.Range("L80").Value = 2439.5557
$ws.Range("M80").Value = -1542.682
$ws.Range("N80").Value = -4435.5557

# Hunk 23: sheet GSM, row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2511.3225
$ws.Range("I83").Value = 2540.682
$ws.Range("J83").Value = 2439.5557
$ws.Range("K83").Value = 12703.41
$ws.Range("L83").Value = 12197.7785
$ws.Range("M83").Value = -7711.41
$ws.Range("N83").Value = -22181.7785

# Hunk 24: sheet WVR, row 2
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 189857.14
$ws.Range("J2").Value = 171500
$ws.Range("L2").Value = 171500
$ws.Range("N2").Value = -171724

# Hunk 25: sheet WVR, row 93
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 29000
$ws.Range("J93").Value = 29000
$ws.Range("L93").Value = 29000
$ws.Range("N93").Value = -33992

# Hunk 26: sheet WVR, row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3543.8462
$ws.Range("I136").Value = 4189.2856
$ws.Range("J136").Value = 2790.8333
$ws.Range("K136").Value = 12567.8568
$ws.Range("L136").Value = 8372.499899999999
$ws.Range("M136").Value = -10017.8568
$ws.Range("N136").Value = -13472.4999
